$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "93.471.35"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +0.43%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.421.20"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.41%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.47"

$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.51%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "620.76"

$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.85%  "

$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +2.74%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.395"

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +1.77%  "

$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.978"

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +5.09%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "3.421.56"

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +1.40%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "43.13"

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +7.30%  "

$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.45%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.32"

$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +4.99%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.287.95"

$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.37%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.059.60"

$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +1.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000247"

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +1.94%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "8.27"

$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +3.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.420.75"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.15"

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +7.80%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.66"

$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +6.49%  "

$ws.Range("B22").Value = "Stellar"

$ws.Range("C22").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.495"

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +8.97%  "

$ws.Range("B23").Value = "SuiNetwork"

$ws.Range("C23").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.40"

$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +8.91%  "

$ws.Range("B24").Value = "BitcoinCash"

$ws.Range("C24").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "500.28"

$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.93%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "6.64"

$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +6.23%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000185"

$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.17%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "95.52"

$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +6.53%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "12.09"

$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +5.45%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.609.89"

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +1.98%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "11.39"

$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +1.19%  "

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.04%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.140"

$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +7.04%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.75"

$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +4.70%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.03"

$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +3.17%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.175"

$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +3.41%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.555"

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +4.77%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "29.17"

$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +2.35%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "568.53"

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.44%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "7.50"

$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.01%  "

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +1.82%  "

$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.02%  "

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.80%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.904"

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.55%  "

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +4.59%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "23.68"

$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.39%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.70"

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.07%  "

$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +5.64%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.50"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.97%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.45"

$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.83%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.13"

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.03%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.17"

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +5.30%  "
